# Update "想去人数" (interest count) values in column F across sheets
# 展览 (Exhibition), 演出 (Performance), 全部类型 (All types)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1095
$ws1.Range("F6").Value  = 613
$ws1.Range("F7").Value  = 584
$ws1.Range("F8").Value  = 1510
$ws1.Range("F10").Value = 1412
$ws1.Range("F11").Value = 3050
$ws1.Range("F12").Value = 553
$ws1.Range("F13").Value = 1718
$ws1.Range("F14").Value = 1779
$ws1.Range("F15").Value = 829
$ws1.Range("F17").Value = 1440
$ws1.Range("F18").Value = 277
$ws1.Range("F19").Value = 70
$ws1.Range("F20").Value = 1174
$ws1.Range("F21").Value = 385
$ws1.Range("F22").Value = 428
$ws1.Range("F23").Value = 53
$ws1.Range("F24").Value = 4629
$ws1.Range("F25").Value = 728
$ws1.Range("F27").Value = 1613
$ws1.Range("F28").Value = 26
$ws1.Range("F29").Value = 76

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 40
$ws2.Range("F14").Value = 21

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 40
$ws4.Range("F15").Value = 1095
$ws4.Range("F17").Value = 613
$ws4.Range("F18").Value = 584
$ws4.Range("F19").Value = 1510
$ws4.Range("F21").Value = 1412
$ws4.Range("F22").Value = 3050
$ws4.Range("F23").Value = 553
$ws4.Range("F24").Value = 1718
$ws4.Range("F25").Value = 1779
$ws4.Range("F26").Value = 829
$ws4.Range("F28").Value = 1440
$ws4.Range("F29").Value = 277
$ws4.Range("F30").Value = 70
$ws4.Range("F33").Value = 1174
$ws4.Range("F34").Value = 385
$ws4.Range("F35").Value = 428
$ws4.Range("F36").Value = 53
$ws4.Range("F37").Value = 4629
$ws4.Range("F38").Value = 728
$ws4.Range("F40").Value = 1613
$ws4.Range("F42").Value = 21
$ws4.Range("F43").Value = 26
$ws4.Range("F44").Value = 76
